$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Q-column attendance marks (value 1) for all student rows except the
# withdrawn-student row (16), mirroring the pattern already present in
# columns D..P for this program week.
$qRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,17,18,19,20)
foreach ($r in $qRows) {
    $ws.Cells.Item($r, 17).Value = 1
}

# New comment on Q2 (column header row) documenting the scholarship award
# ceremony participation.
$ws.Range("Q2").AddComment("장학수여식 참여`n해교담임선생님께 문의") | Out-Null

# Update the last-selected cell to match the author's final selection.
$ws.Range("Q27").Select()
